$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on the Price/Volume columns so numeric-looking
# strings (e.g. "1.001", "10.60") are preserved exactly as text,
# matching the original inlineStr cell type.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range('D2').Value = '27.667.26'
$ws.Range('E2').Value = '  -0.87%  '
$ws.Range('D3').Value = '1.893.73'
$ws.Range('E3').Value = '  +1.04%  '
$ws.Range('D4').Value = '1.001'
$ws.Range('E4').Value = '  -1.01%  '
$ws.Range('D5').Value = '312.84'
$ws.Range('E5').Value = '  -0.31%  '
$ws.Range('E6').Value = '  -0.98%  '
$ws.Range('D7').Value = '0.4858'
$ws.Range('E7').Value = '  +0.57%  '
$ws.Range('E8').Value = '  -0.85%  '
$ws.Range('D9').Value = '0.07333'
$ws.Range('E9').Value = '  -0.70%  '
$ws.Range('D10').Value = '0.9156'
$ws.Range('E10').Value = '  -2.55%  '
$ws.Range('D11').Value = '20.53'
$ws.Range('E11').Value = '  -2.44%  '
$ws.Range('D12').Value = '0.07697'
$ws.Range('E12').Value = '  -1.51%  '
$ws.Range('D13').Value = '1.894.84'
$ws.Range('E13').Value = '  -0.24%  '
$ws.Range('D14').Value = '5.478'
$ws.Range('D15').Value = '6.604'
$ws.Range('E15').Value = '  -0.23%  '
$ws.Range('E16').Value = '  -0.02%  '
$ws.Range('E17').Value = '  -0.99%  '
$ws.Range('D18').Value = '0.000008799'
$ws.Range('E18').Value = '  -1.05%  '
$ws.Range('E19').Value = '  -0.97%  '
$ws.Range('D20').Value = '27.703.43'
$ws.Range('E20').Value = '  -0.82%  '
$ws.Range('E21').Value = '  -2.38%  '
$ws.Range('D22').Value = '5.127'
$ws.Range('E22').Value = '  +0.02%  '
$ws.Range('D23').Value = '2.095.63'
$ws.Range('E23').Value = '  -1.33%  '
$ws.Range('D24').Value = '10.74'
$ws.Range('E24').Value = '  -0.94%  '
$ws.Range('D25').Value = '1.902'
$ws.Range('E25').Value = '  -2.50%  '
$ws.Range('D26').Value = '153.51'
$ws.Range('E26').Value = '  -2.01%  '
$ws.Range('D27').Value = '18.36'
$ws.Range('E27').Value = '  -1.19%  '
$ws.Range('D28').Value = '2.146'
$ws.Range('E28').Value = '  +4.20%  '
$ws.Range('D29').Value = '115.77'
$ws.Range('D30').Value = '4.912'
$ws.Range('E30').Value = '  -1.53%  '
$ws.Range('D31').Value = '0.08917'
$ws.Range('E31').Value = '  -0.03%  '
$ws.Range('D32').Value = '3.184'
$ws.Range('E32').Value = '  -4.45%  '
$ws.Range('D33').Value = '1.221'
$ws.Range('E33').Value = '  -0.27%  '
$ws.Range('D34').Value = '0.7658'
$ws.Range('E34').Value = '  -0.02%  '
$ws.Range('D35').Value = '4.643'
$ws.Range('E35').Value = '  -0.61%  '
$ws.Range('B36').Value = 'VeChain'
$ws.Range('C36').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D36').Value = '0.02041'
$ws.Range('E36').Value = '  -0.33%  '
$ws.Range('B37').Value = 'RenderToken'
$ws.Range('C37').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D37').Value = '2.524'
$ws.Range('E37').Value = '  -7.26%  '
$ws.Range('B38').Value = 'TrustWalletToken'
$ws.Range('C38').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D38').Value = '1.094'
$ws.Range('E38').Value = '  -3.57%  '
$ws.Range('B39').Value = 'Hedera'
$ws.Range('C39').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D39').Value = '0.05275'
$ws.Range('E39').Value = '  -1.82%  '
$ws.Range('B40').Value = 'TheSandbox'
$ws.Range('C40').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D40').Value = '0.5475'
$ws.Range('E40').Value = '  -3.18%  '
$ws.Range('B41').Value = 'MXToken'
$ws.Range('C41').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D41').Value = '2.981'
$ws.Range('E41').Value = '  -0.49%  '
$ws.Range('B42').Value = 'FraxShare'
$ws.Range('C42').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D42').Value = '6.912'
$ws.Range('E42').Value = '  -2.21%  '
$ws.Range('B43').Value = 'Aptos'
$ws.Range('C43').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D43').Value = '8.475'
$ws.Range('E43').Value = '  -1.31%  '
$ws.Range('D44').Value = '0.1516'
$ws.Range('E44').Value = '  -1.13%  '
$ws.Range('B45').Value = 'Quant'
$ws.Range('C45').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D45').Value = '110.94'
$ws.Range('E45').Value = '  +5.29%  '
$ws.Range('B46').Value = 'EnergySwap'
$ws.Range('C46').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D46').Value = '10.60'
$ws.Range('E46').Value = '  -1.51%  '
$ws.Range('B47').Value = 'Decentraland'
$ws.Range('C47').Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range('D47').Value = '0.4795'
$ws.Range('E47').Value = '  -2.21%  '
$ws.Range('B48').Value = 'PaxDollar'
$ws.Range('C48').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range('D48').Value = '1.001'
$ws.Range('E48').Value = '  -1.08%  '
$ws.Range('B49').Value = 'NEARProtocol'
$ws.Range('C49').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D49').Value = '1.636'
$ws.Range('E49').Value = '  -2.42%  '
$ws.Range('B50').Value = 'Aave'
$ws.Range('C50').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D50').Value = '67.37'
$ws.Range('E50').Value = '  -0.63%  '
$ws.Range('B51').Value = 'Cronos'
$ws.Range('C51').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D51').Value = '0.06051'
$ws.Range('E51').Value = '  -0.97%  '

# Restore default cell style (keeps cells styleless, as in the original)
$ws.Range("D2:E51").Style = "Normal"
